# Atualização de bases das ligas, do dia: 12-06-2024 às 23:38
# Swap the full data (columns B:AD) between each pair of rows listed below,
# leaving column A (the sequential id) untouched in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(107, 108),
    @(128, 129),
    @(143, 145),
    @(148, 149)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B${r1}:AD${r1}")
    $range2 = $ws.Range("B${r2}:AD${r2}")

    $vals1 = $range1.Value()
    $vals2 = $range2.Value()

    $range1.Value = $vals2
    $range2.Value = $vals1
}
